$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the saved window position (matches the author's recorded window move)
$wb.Windows.Item(1).Left = 620

# Add a defined name that references a table
$wb.Names.Add("named_reference_to_table", "=Table1[#All]")

# Add labels showing the name and its formula text in F7:G7
$ws.Range("F7").Value = "named_reference_to_table"

# G7 holds the literal text "=Table1[#All]" (not an actual formula). Enter it
# as a formula that evaluates to the desired string, then convert the cell to
# its value in place so it lands in the sheet as plain text.
$ws.Range("G7").Formula = '="=Table1[#All]"'
$ws.Range("G7").Copy()
$ws.Range("G7").PasteSpecial(-4163)

# Update selection to G7
$ws.Range("G7").Select()
